$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-17 16:42:18"
# Target stored column width is 17.2159881591797 characters; the ColumnWidth
# COM property quantizes to whole-pixel steps (1/6 character), so 16.33 is the
# closest input that lands on the nearest achievable stored width (17.1667).
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-17 16:41:55"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-17 16:42:18"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
